$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these price cells as text (they look numeric but are stored as text,
# matching the source data, e.g. "208.57" not 208.57).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "25.547.20"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.588.56"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "208.57"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  -4.75%  "
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "17.73"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").Value = "0.0783"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "1.586.81"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "0.507"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "25.588.14"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "60.31"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "0.0₃0711"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "187.64"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "9.32"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").Value = "5.92"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("D26").Value = "140.53"
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("D27").Value = "1.69"
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("D28").Value = "6.49"
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("D29").Value = "14.92"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("D31").Value = "0.0465"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").Value = "3.06"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  -4.36%  "
$ws.Range("D34").Value = "2.40"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "1.46"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").Value = "1.088.03"
$ws.Range("E36").Value = "  -4.16%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").Value = "0.0150"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").Value = "0.778"
$ws.Range("E40").Value = "  -10.11%  "
$ws.Range("D41").Value = "0.494"
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("D42").Value = "95.19"
$ws.Range("E42").Value = "  -3.28%  "
$ws.Range("D43").Value = "1.724.69"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("D45").Value = "0.733"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("E46").Value = "  -7.78%  "
$ws.Range("D47").Value = "52.88"
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("D48").Value = "0.0508"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("D49").Value = "1.42"
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("D50").Value = "0.408"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("E51").Value = "  -0.24%  "
